# #249: added test case and updated test file
# Adds a small 4-column Excel Table (ListObject) "Table1" to the right of the
# existing data (O4:R13), with header row Column1..Column4, autofiltered and
# styled with TableStyleMedium2. Also widens the new columns and moves the
# active selection, matching the updated NvPr.xlsx test workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new table (the ListObjects.Add call below also writes
# these values/headers, but setting them explicitly keeps things clear).
$ws.Range("O4").Value = "Column1"
$ws.Range("P4").Value = "Column2"
$ws.Range("Q4").Value = "Column3"
$ws.Range("R4").Value = "Column4"

# Create the table over O4:R13 using the existing header row.
$lo = $ws.ListObjects.Add(1, $ws.Range("O4:R13"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleMedium2"
$lo.Comment = "Luke, I am your father... seriously..."

# Widen the new columns (O:R) to match the authored workbook.
$ws.Range("O4:R4").ColumnWidth = 10.166666666666666

# Move the active selection as in the updated workbook.
$ws.Range("R24").Select() | Out-Null
